$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting rows 40:122 down to 41:123
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record's data
$ws.Cells.Item(40, 1).Value = 11
$ws.Cells.Item(40, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(40, 3).Value = "Bíobío"
$ws.Cells.Item(40, 4).Value = 44544
$ws.Cells.Item(40, 5).Value = 8
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100108
$ws.Cells.Item(40, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(40, 9).Value = 100108005
$ws.Cells.Item(40, 10).Value = "Piña"
$ws.Cells.Item(40, 11).Value = "Caramelo"
$ws.Cells.Item(40, 12).Value = "Segunda"
$ws.Cells.Item(40, 13).Value = 200
$ws.Cells.Item(40, 14).Value = 16000
$ws.Cells.Item(40, 15).Value = 17000
$ws.Cells.Item(40, 16).Value = 16500
$ws.Cells.Item(40, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(40, 18).Value = "Ecuador"
$ws.Cells.Item(40, 19).Value = 1179
$ws.Cells.Item(40, 20).Value = 14
